$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 182. This shifts all existing rows
# 182-290 down to 183-291 (preserving their data), and creates a
# new, empty row 182 ready to be populated with the new record.
$ws.Rows.Item(182).Insert()

# Populate the newly inserted row 182 with the new data record.
$ws.Cells.Item(182, 1).Value = 10
$ws.Cells.Item(182, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(182, 3).Value = "La Araucanía"
$ws.Cells.Item(182, 4).Value = 44719
$ws.Cells.Item(182, 5).Value = 9
$ws.Cells.Item(182, 6).Value = 100112001
$ws.Cells.Item(182, 7).Value = "Berenjena"
$ws.Cells.Item(182, 8).Value = "Sin especificar"
$ws.Cells.Item(182, 9).Value = "Primera"
$ws.Cells.Item(182, 10).Value = 50
$ws.Cells.Item(182, 11).Value = 10000
$ws.Cells.Item(182, 12).Value = 10000
$ws.Cells.Item(182, 13).Value = 10000
$ws.Cells.Item(182, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(182, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(182, 16).Value = 167
$ws.Cells.Item(182, 17).Value = 60
$ws.Cells.Item(182, 18).Value = "Hortaliza"
